# New weekly price observation (Espinaca, Vega Monumental Concepción) is
# inserted as a new row at position 28, pushing the existing rows 28-82
# down to rows 29-83 (mirrors a new week being added to the top of the
# chronological price history kept in this sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(28).Insert()

$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = 44757
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = 100112012
$ws.Cells.Item(28, 7).Value = "Espinaca"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 40
$ws.Cells.Item(28, 11).Value = 8000
$ws.Cells.Item(28, 12).Value = 8500
$ws.Cells.Item(28, 13).Value = 8250
$ws.Cells.Item(28, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 825
$ws.Cells.Item(28, 17).Value = 10
$ws.Cells.Item(28, 18).Value = "Hortaliza"

Write-Host "done"
